# Update cryptos list values (Price and Volume(1h) columns) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.884.90"
$ws.Range("E2").Value = "  +4.87%  "
$ws.Range("D3").Value = "2.268.70"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'302.08"
$ws.Range("E5").Value = "  +3.72%  "
$ws.Range("D6").Value = "'92.14"
$ws.Range("E6").Value = "  +6.32%  "
$ws.Range("D7").Value = "'0.533"
$ws.Range("E7").Value = "  +3.70%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +3.44%  "
$ws.Range("E10").Value = "  +8.75%  "
$ws.Range("D11").Value = "'32.25"
$ws.Range("E12").Value = "  +2.66%  "
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("E14").Value = "  +3.65%  "
$ws.Range("D15").Value = "2.619.88"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").Value = "'14.14"
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("D17").Value = "2.263.44"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("D19").Value = "41.791.82"
$ws.Range("E19").Value = "  +4.81%  "
$ws.Range("D20").Value = "'12.08"
$ws.Range("E20").Value = "  +9.04%  "
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").Value = "'5.93"
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("D23").Value = "'66.97"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("D24").Value = "'240.67"
$ws.Range("D25").Value = "'2.55"
$ws.Range("E25").Value = "  +3.82%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  +3.93%  "
$ws.Range("E28").Value = "  +3.98%  "
$ws.Range("E29").Value = "  +6.67%  "
$ws.Range("D30").Value = "'9.62"
$ws.Range("E30").Value = "  +4.26%  "
$ws.Range("D31").Value = "'159.14"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").Value = "'33.82"
$ws.Range("E32").Value = "  +6.27%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'5.15"
$ws.Range("E34").Value = "  +3.99%  "
$ws.Range("D35").Value = "'0.0745"
$ws.Range("E35").Value = "  +4.56%  "
$ws.Range("E36").Value = "  +3.10%  "
$ws.Range("E38").Value = "  +5.09%  "
$ws.Range("D39").Value = "'0.116"
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("D40").Value = "'16.51"
$ws.Range("E40").Value = "  +8.48%  "
$ws.Range("D41").Value = "'1.80"
$ws.Range("E41").Value = "  +3.94%  "
$ws.Range("E42").Value = "  +6.10%  "
$ws.Range("D43").Value = "2.065.35"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D44").Value = "'19.53"
$ws.Range("E44").Value = "  +8.89%  "
$ws.Range("E45").Value = "  +2.64%  "
$ws.Range("D46").Value = "'10.08"
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("E47").Value = "  +6.38%  "
$ws.Range("D48").Value = "'2.03"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("E49").Value = "  +3.83%  "
$ws.Range("E50").Value = "  +3.19%  "
$ws.Range("D51").Value = "'51.59"
$ws.Range("E51").Value = "  +5.27%  "
